$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A70").Value = "2025-04-29 10:11:40"
$ws.Range("B70").Value = 203
